$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A21: was stored as inline string "71277620", should become numeric 71277620
$ws.Range("A21").Value = 71277620

# Add new row 22: redemption of 100 points for phone 79174445 at given timestamp.
# The phone number must stay stored as text (matches the other "phone" entries
# in the sheet), so temporarily force a text number format before assigning the
# value - otherwise the numeric-looking string gets auto-coerced to a number -
# then clear the formatting again so the new row's cells end up unstyled, same
# as their neighbours.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "79174445"
$ws.Range("A22").ClearFormats()
$ws.Range("B22").Value = 100
$ws.Range("C22").Value = "2025-08-18T17:42:29"
